# Updates cryptos price list values to match latest scrape (coinranking.com)
# Source diff: commit "Updated cryptos list on Tue Mar  7 07:54:19 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Cells.Item(2, 4).Value = '22.438.15'
$ws.Cells.Item(2, 5).Value = '  +0.12%  '

# Row 3 (Ethereum)
$ws.Cells.Item(3, 4).Value = '1.573.34'
$ws.Cells.Item(3, 5).Value = '  +0.65%  '

# Row 4 (TetherUSD)
$ws.Cells.Item(4, 5).Value = '  +0.05%  '

# Row 6 (BNB)
$ws.Cells.Item(6, 4).Value = '''287.88'

# Row 7 (XRP)
$ws.Cells.Item(7, 4).Value = '''0.3708'
$ws.Cells.Item(7, 5).Value = '  +2.04%  '

# Row 8 (OKB)
$ws.Cells.Item(8, 4).Value = '''47.34'
$ws.Cells.Item(8, 5).Value = '  -2.55%  '

# Row 9 (Cardano)
$ws.Cells.Item(9, 4).Value = '''0.3315'
$ws.Cells.Item(9, 5).Value = '  -0.63%  '

# Row 10 (Polygon)
$ws.Cells.Item(10, 5).Value = '  +2.52%  '

# Row 11 (Dogecoin)
$ws.Cells.Item(11, 4).Value = '''0.07503'
$ws.Cells.Item(11, 5).Value = '  +1.26%  '

# Row 12 (BinanceUSD)
$ws.Cells.Item(12, 5).Value = '  +0.09%  '

# Row 13 (Solana)
$ws.Cells.Item(13, 4).Value = '''20.73'
$ws.Cells.Item(13, 5).Value = '  -0.37%  '

# Row 14 (Polkadot)
$ws.Cells.Item(14, 4).Value = '''5.932'

# Row 15 (Chainlink)
$ws.Cells.Item(15, 4).Value = '''6.914'
$ws.Cells.Item(15, 5).Value = '  +0.28%  '

# Row 16 (WrappedEther)
$ws.Cells.Item(16, 4).Value = '1.563.36'
$ws.Cells.Item(16, 5).Value = '  -0.17%  '

# Row 17 (ShibaInu)
$ws.Cells.Item(17, 4).Value = '''0.00001114'
$ws.Cells.Item(17, 5).Value = '  +0.90%  '

# Row 18 (Litecoin)
$ws.Cells.Item(18, 4).Value = '''88.33'
$ws.Cells.Item(18, 5).Value = '  +0.14%  '

# Row 19 (TRON)
$ws.Cells.Item(19, 4).Value = '''0.06728'
$ws.Cells.Item(19, 5).Value = '  +0.55%  '

# Row 20 (Dai)
$ws.Cells.Item(20, 2).Value = 'Uniswap'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(20, 4).Value = '''6.407'
$ws.Cells.Item(20, 5).Value = '  +0.84%  '

# Row 21 (Uniswap)
$ws.Cells.Item(21, 2).Value = 'Dai'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(21, 4).Value = '''1.000'
$ws.Cells.Item(21, 5).Value = '  -0.01%  '

# Row 22 (Avalanche)
$ws.Cells.Item(22, 4).Value = '''16.49'

# Row 23 (Cosmos)
$ws.Cells.Item(23, 4).Value = '''11.97'
$ws.Cells.Item(23, 5).Value = '  -0.22%  '

# Row 24 (WrappedBTC)
$ws.Cells.Item(24, 4).Value = '22.424.67'

# Row 25 (Toncoin)
$ws.Cells.Item(25, 4).Value = '''2.371'
$ws.Cells.Item(25, 5).Value = '  -1.87%  '

# Row 26 (LidoDAOToken)
$ws.Cells.Item(26, 4).Value = '''2.624'
$ws.Cells.Item(26, 5).Value = '  +2.58%  '

# Row 27 (Monero)
$ws.Cells.Item(27, 4).Value = '''150.79'
$ws.Cells.Item(27, 5).Value = '  +0.70%  '

# Row 28 (EthereumClassic)
$ws.Cells.Item(28, 4).Value = '''19.55'
$ws.Cells.Item(28, 5).Value = '  +0.74%  '

# Row 29 (HuobiToken)
$ws.Cells.Item(29, 4).Value = '''4.947'
$ws.Cells.Item(29, 5).Value = '  -0.94%  '

# Row 30 (BitcoinCash)
$ws.Cells.Item(30, 4).Value = '''124.95'
$ws.Cells.Item(30, 5).Value = '  +1.43%  '

# Row 31 (WrappedliquidstakedEther2.0)
$ws.Cells.Item(31, 4).Value = '1.739.82'
$ws.Cells.Item(31, 5).Value = '  +0.01%  '

# Row 32 (ImmutableX)
$ws.Cells.Item(32, 4).Value = '''1.090'
$ws.Cells.Item(32, 5).Value = '  +2.81%  '

# Row 33 (Filecoin)
$ws.Cells.Item(33, 4).Value = '''6.079'
$ws.Cells.Item(33, 5).Value = '  -0.95%  '

# Row 34 (WEMIXTOKEN)
$ws.Cells.Item(34, 5).Value = '  -0.56%  '

# Row 35 (FraxShare)
$ws.Cells.Item(35, 4).Value = '''9.864'
$ws.Cells.Item(35, 5).Value = '  +0.51%  '

# Row 36 (Stellar)
$ws.Cells.Item(36, 4).Value = '''0.08321'
$ws.Cells.Item(36, 5).Value = '  +0.96%  '

# Row 37 (VeChain)
$ws.Cells.Item(37, 4).Value = '''0.02441'
$ws.Cells.Item(37, 5).Value = '  +1.76%  '

# Row 38 (TrustWalletToken)
$ws.Cells.Item(38, 5).Value = '  -0.19%  '

# Row 39 (Hedera)
$ws.Cells.Item(39, 4).Value = '''0.06378'

# Row 40 (Algorand)
$ws.Cells.Item(40, 4).Value = '''0.2214'
$ws.Cells.Item(40, 5).Value = '  +0.22%  '

# Row 41 (InternetComputer(DFINITY))
$ws.Cells.Item(41, 4).Value = '''5.327'
$ws.Cells.Item(41, 5).Value = '  -0.09%  '

# Row 42 (Aptos)
$ws.Cells.Item(42, 5).Value = '  +2.00%  '

# Row 43 (TheSandbox)
$ws.Cells.Item(43, 4).Value = '''0.6231'
$ws.Cells.Item(43, 5).Value = '  +2.41%  '

# Row 44 (EnergySwap)
$ws.Cells.Item(44, 2).Value = 'Frax'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(44, 4).Value = '''1.000'
$ws.Cells.Item(44, 5).Value = '  +0.05%  '

# Row 45 (Decentraland)
$ws.Cells.Item(45, 2).Value = 'EnergySwap'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(45, 4).Value = '''13.96'
$ws.Cells.Item(45, 5).Value = '  +0.97%  '

# Row 46 (PancakeSwap)
$ws.Cells.Item(46, 2).Value = 'Decentraland'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(46, 4).Value = '''0.6038'
$ws.Cells.Item(46, 5).Value = '  +4.86%  '

# Row 47 (NEARProtocol)
$ws.Cells.Item(47, 2).Value = 'PancakeSwap'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(47, 4).Value = '''3.773'
$ws.Cells.Item(47, 5).Value = '  +0.46%  '

# Row 48 (Quant)
$ws.Cells.Item(48, 2).Value = 'NEARProtocol'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(48, 4).Value = '''2.037'
$ws.Cells.Item(48, 5).Value = '  +1.22%  '

# Row 49 (EOS)
$ws.Cells.Item(49, 2).Value = 'Quant'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(49, 4).Value = '''124.77'
$ws.Cells.Item(49, 5).Value = '  +0.21%  '

# Row 50 (Cronos)
$ws.Cells.Item(50, 2).Value = 'EOS'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Cells.Item(50, 4).Value = '''1.205'
$ws.Cells.Item(50, 5).Value = '  -0.75%  '

# Row 51 (Aave)
$ws.Cells.Item(51, 2).Value = 'Cronos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(51, 4).Value = '''0.07193'
$ws.Cells.Item(51, 5).Value = '  -0.22%  '
